$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column for all existing data rows (2-409)
#    from 45192 to 45202, keeping the existing date style/format.
$ws.Range("C2:C409").Value = 45202

# 2. Append two new report rows (410 and 411) with the same column layout
#    as the existing rows (A-E, G-Q, and an empty wrap-formatted R cell).

# --- Row 410 ---
$ws.Range("A410").Value = "A 46583-2023"
$ws.Range("B410").Value = 45198
$ws.Range("C410").Value = 45202
$ws.Range("B410:C410").NumberFormat = "YYYY-MM-DD"
$ws.Range("D410").Value = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E410").Value = "YDRE"
$ws.Range("G410").Value = 10.9
$ws.Range("H410").Value = 0
$ws.Range("I410").Value = 0
$ws.Range("J410").Value = 0
$ws.Range("K410").Value = 0
$ws.Range("L410").Value = 0
$ws.Range("M410").Value = 0
$ws.Range("N410").Value = 0
$ws.Range("O410").Value = 0
$ws.Range("P410").Value = 0
$ws.Range("Q410").Value = 0
$ws.Range("R410").Value = ""
$ws.Range("R410").WrapText = $true

# --- Row 411 ---
$ws.Range("A411").Value = "A 46578-2023"
$ws.Range("B411").Value = 45198
$ws.Range("C411").Value = 45202
$ws.Range("B411:C411").NumberFormat = "YYYY-MM-DD"
$ws.Range("D411").Value = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E411").Value = "YDRE"
$ws.Range("G411").Value = 1.9
$ws.Range("H411").Value = 0
$ws.Range("I411").Value = 0
$ws.Range("J411").Value = 0
$ws.Range("K411").Value = 0
$ws.Range("L411").Value = 0
$ws.Range("M411").Value = 0
$ws.Range("N411").Value = 0
$ws.Range("O411").Value = 0
$ws.Range("P411").Value = 0
$ws.Range("Q411").Value = 0
$ws.Range("R411").Value = ""
$ws.Range("R411").WrapText = $true

# Match the source workbook's row-height stamping: every row gains an
# explicit custom height of 15 once a new row is appended below it,
# except for the very last row in the sheet.
$ws.Rows.Item(409).RowHeight = 15
$ws.Rows.Item(410).RowHeight = 15
